$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.963.93'
$ws.Range("E2").Value = '  -0.59%  '
$ws.Range("D3").Value = '1.856.63'
$ws.Range("E3").Value = '  -1.09%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '312.24'
$ws.Range("E5").Value = '  -0.54%  '
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("E7").Value = '  +1.46%  '
$ws.Range("D8").Value = '0.3827'
$ws.Range("E8").Value = '  -0.39%  '
$ws.Range("D9").Value = '0.08217'
$ws.Range("E9").Value = '  -4.90%  '
$ws.Range("D10").Value = '1.108'
$ws.Range("E10").Value = '  -0.88%  '
$ws.Range("D11").Value = '41.46'
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("D12").Value = '6.175'
$ws.Range("E12").Value = '  -2.41%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.866.79'
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("B14").Value = 'Solana'
$ws.Range("C14").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D14").Value = '20.47'
$ws.Range("E14").Value = '  -1.00%  '
$ws.Range("D15").Value = '7.247'
$ws.Range("E15").Value = '  +0.93%  '
$ws.Range("D16").Value = '1.005'
$ws.Range("D17").Value = '0.00001096'
$ws.Range("E17").Value = '  -0.59%  '
$ws.Range("D18").Value = '90.42'
$ws.Range("E18").Value = '  -0.71%  '
$ws.Range("D19").Value = '0.06647'
$ws.Range("E19").Value = '  +0.36%  '
$ws.Range("D20").Value = '17.64'
$ws.Range("E20").Value = '  -2.97%  '
$ws.Range("E21").Value = '  +0.13%  '
$ws.Range("D22").Value = '6.004'
$ws.Range("E22").Value = '  -1.58%  '
$ws.Range("D23").Value = '28.006.75'
$ws.Range("E23").Value = '  -0.54%  '
$ws.Range("D24").Value = '11.03'
$ws.Range("E24").Value = '  -3.50%  '
$ws.Range("D25").Value = '2.257'
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("D26").Value = '2.072.82'
$ws.Range("E26").Value = '  -1.30%  '
$ws.Range("D27").Value = '2.501'
$ws.Range("E27").Value = '  -1.85%  '
$ws.Range("D28").Value = '157.07'
$ws.Range("E28").Value = '  -0.28%  '
$ws.Range("D29").Value = '20.41'
$ws.Range("E29").Value = '  -1.46%  '
$ws.Range("D30").Value = '124.53'
$ws.Range("E30").Value = '  -1.92%  '
$ws.Range("D31").Value = '0.1064'
$ws.Range("E31").Value = '  +1.17%  '
$ws.Range("D32").Value = '1.030'
$ws.Range("E32").Value = '  -2.81%  '
$ws.Range("D33").Value = '5.899'
$ws.Range("E33").Value = '  +5.48%  '
$ws.Range("E34").Value = '  -0.45%  '
$ws.Range("D35").Value = '9.376'
$ws.Range("E35").Value = '  -2.99%  '
$ws.Range("D36").Value = '0.06508'
$ws.Range("E36").Value = '  -0.52%  '
$ws.Range("D37").Value = '0.02409'
$ws.Range("E37").Value = '  -1.64%  '
$ws.Range("D38").Value = '0.2174'
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("D39").Value = '0.6527'
$ws.Range("E39").Value = '  +2.30%  '
$ws.Range("D40").Value = '1.191'
$ws.Range("E40").Value = '  -1.03%  '
$ws.Range("D41").Value = '4.970'
$ws.Range("E41").Value = '  +1.49%  '
$ws.Range("D42").Value = '1.210'
$ws.Range("E42").Value = '  -2.42%  '
$ws.Range("D43").Value = '11.14'
$ws.Range("E43").Value = '  -3.53%  '
$ws.Range("D44").Value = '0.6139'
$ws.Range("E44").Value = '  +2.56%  '
$ws.Range("D45").Value = '13.01'
$ws.Range("E45").Value = '  -1.66%  '
$ws.Range("D46").Value = '1.278'
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").Value = '3.673'
$ws.Range("E47").Value = '  -0.09%  '
$ws.Range("D48").Value = '2.009'
$ws.Range("E48").Value = '  +0.93%  '
$ws.Range("D49").Value = '1.211'
$ws.Range("E49").Value = '  -1.65%  '
$ws.Range("D50").Value = '120.79'
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("D51").Value = '78.14'
$ws.Range("E51").Value = '  -2.15%  '
